$wb = $excel.ActiveWorkbook

# --- Sheet: Significant Components ---
$ws1 = $wb.Worksheets.Item("Significant Components")
$ws1.Range("C2").Value = "['QEXTRCT' 'QESL' 'QHISPC' 'QEDLESHI' 'PPUNIT' 'QNOHLTH' 'QSERV' 'QFHH'`n 'PERCAP']"
$ws1.Range("C3").Value = "['QNOAUTO' 'QPOVTY' 'MDGRENT' 'QRENTER' 'QFAM']"
$ws1.Range("C4").Value = "['QRICH' 'PERCAP' 'MDHSEVAL']"
$ws1.Range("C5").Value = "['QRENTER' 'QAGEDEP' 'QSSBEN' 'MEDAGE']"
$ws1.Range("C6").Value = "['QAGEDEP' 'QFEMLBR' 'QFEMALE']"

# --- Sheet: Loading Factors ---
$ws2 = $wb.Worksheets.Item("Loading Factors")
$ws2.Range("A2").Value = "QEXTRCT"
$ws2.Range("B2").Value = 0.6858736575524721
$ws2.Range("C2").Value = 0.1686737909858408
$ws2.Range("D2").Value = 0.1678446216009458
$ws2.Range("E2").Value = 0.06581878842978414
$ws2.Range("F2").Value = -0.3801092957356946
$ws2.Range("A3").Value = "QESL"
$ws2.Range("B3").Value = 0.7871772107263549
$ws2.Range("C3").Value = 0.2433220629746484
$ws2.Range("D3").Value = 0.1589091485125755
$ws2.Range("E3").Value = -0.03444785117753878
$ws2.Range("F3").Value = -0.306829493692647
$ws2.Range("B4").Value = 0.8256792357545928
$ws2.Range("C4").Value = 0.2303082105529956
$ws2.Range("D4").Value = 0.3322223344528948
$ws2.Range("E4").Value = -0.08286797780992657
$ws2.Range("F4").Value = -0.1013539821395284
$ws2.Range("A5").Value = "QEDLESHI"
$ws2.Range("B5").Value = 0.8535540399647311
$ws2.Range("C5").Value = 0.2942599856188969
$ws2.Range("D5").Value = 0.2278342548178292
$ws2.Range("E5").Value = 0.03093405418756537
$ws2.Range("F5").Value = -0.1792954128540681
$ws2.Range("A6").Value = "PPUNIT"
$ws2.Range("B6").Value = 0.8390350525801112
$ws2.Range("C6").Value = -0.3911401996634874
$ws2.Range("D6").Value = 0.02055685611992847
$ws2.Range("E6").Value = -0.104979825283385
$ws2.Range("F6").Value = 0.06282084972806128
$ws2.Range("A7").Value = "QNOHLTH"
$ws2.Range("B7").Value = 0.5868900526934003
$ws2.Range("C7").Value = 0.3769984876564605
$ws2.Range("D7").Value = 0.4345518397036272
$ws2.Range("E7").Value = -0.112250758320281
$ws2.Range("F7").Value = -0.1988058259654437
$ws2.Range("A8").Value = "QSERV"
$ws2.Range("B8").Value = 0.5235844267086339
$ws2.Range("C8").Value = 0.3723549677993631
$ws2.Range("D8").Value = 0.3418577957046304
$ws2.Range("E8").Value = -0.1627309474076945
$ws2.Range("F8").Value = 0.08706326963143393
$ws2.Range("A9").Value = "QFHH"
$ws2.Range("B9").Value = 0.5811161255868825
$ws2.Range("C9").Value = 0.2311338953133625
$ws2.Range("D9").Value = 0.2397958574662523
$ws2.Range("E9").Value = -0.06720382108618962
$ws2.Range("F9").Value = 0.3153982269937073
$ws2.Range("A10").Value = "QRICH"
$ws2.Range("B10").Value = 0.1482951756887093
$ws2.Range("C10").Value = 0.3316849231928604
$ws2.Range("D10").Value = 0.8441686390780203
$ws2.Range("E10").Value = -0.1720039450954542
$ws2.Range("F10").Value = 0.009111125661337663
$ws2.Range("A11").Value = "PERCAP"
$ws2.Range("B11").Value = 0.4617240831990917
$ws2.Range("C11").Value = 0.2826832742982029
$ws2.Range("D11").Value = 0.7178886796790139
$ws2.Range("E11").Value = -0.2445429701812348
$ws2.Range("F11").Value = 0.03072561909958982
$ws2.Range("A12").Value = "MDHSEVAL"
$ws2.Range("B12").Value = 0.3403125200052169
$ws2.Range("C12").Value = 0.1154671295287299
$ws2.Range("D12").Value = 0.8271739462414933
$ws2.Range("E12").Value = -0.06457556998133435
$ws2.Range("F12").Value = 0.02123961798656674
$ws2.Range("B13").Value = 0.1441323019310419
$ws2.Range("C13").Value = 0.7131384755366393
$ws2.Range("D13").Value = 0.06639925210386499
$ws2.Range("E13").Value = -0.04474492005653481
$ws2.Range("F13").Value = 0.008266724356157294
$ws2.Range("A14").Value = "QPOVTY"
$ws2.Range("B14").Value = 0.4087840692542938
$ws2.Range("C14").Value = 0.5110704235684644
$ws2.Range("D14").Value = 0.188001844083467
$ws2.Range("E14").Value = -0.3313326743756278
$ws2.Range("F14").Value = 0.02863860231725125
$ws2.Range("A15").Value = "MDGRENT"
$ws2.Range("B15").Value = -0.1975754434779801
$ws2.Range("C15").Value = -0.535177679099883
$ws2.Range("D15").Value = -0.3165122544058664
$ws2.Range("E15").Value = 0.0218253183896386
$ws2.Range("F15").Value = 0.1164025804462911
$ws2.Range("A16").Value = "QRENTER"
$ws2.Range("B16").Value = -0.05014761453986846
$ws2.Range("C16").Value = 0.7318178143195329
$ws2.Range("D16").Value = 0.2033261131586872
$ws2.Range("E16").Value = -0.4630674592357272
$ws2.Range("F16").Value = -0.09625206428947632
$ws2.Range("B17").Value = 0.2131259487493579
$ws2.Range("C17").Value = 0.5156558124547915
$ws2.Range("D17").Value = 0.2584929403896893
$ws2.Range("E17").Value = -0.1667465807900852
$ws2.Range("F17").Value = 0.1380807234771666
$ws2.Range("A18").Value = "QAGEDEP"
$ws2.Range("B18").Value = 0.004237817472272597
$ws2.Range("C18").Value = -0.1129305024413651
$ws2.Range("D18").Value = -0.1037482651232172
$ws2.Range("E18").Value = 0.6808559555794382
$ws2.Range("F18").Value = 0.5744249450780913
$ws2.Range("A19").Value = "QSSBEN"
$ws2.Range("B19").Value = 0.04206845476943216
$ws2.Range("C19").Value = -0.09549120345561292
$ws2.Range("D19").Value = -0.05344152197055597
$ws2.Range("E19").Value = 0.7693356006567849
$ws2.Range("F19").Value = 0.123001187274755
$ws2.Range("A20").Value = "MEDAGE"
$ws2.Range("B20").Value = -0.3199370047601924
$ws2.Range("C20").Value = -0.3055027584054265
$ws2.Range("D20").Value = -0.2684234148189786
$ws2.Range("E20").Value = 0.7643515739670108
$ws2.Range("F20").Value = -0.05237667736199582
$ws2.Range("A21").Value = "QFEMLBR"
$ws2.Range("B21").Value = -0.1958049799799623
$ws2.Range("C21").Value = 0.0731579741974871
$ws2.Range("D21").Value = 0.09900758791837097
$ws2.Range("E21").Value = -0.01338488484407241
$ws2.Range("F21").Value = 0.7610846981779612
$ws2.Range("A22").Value = "QFEMALE"
$ws2.Range("B22").Value = -0.06730619307533976
$ws2.Range("C22").Value = -0.0489382739587247
$ws2.Range("D22").Value = -0.04702348251678303
$ws2.Range("E22").Value = 0.2326601383396722
$ws2.Range("F22").Value = 0.7981094768299711

# --- Sheet: All Refactor Variances ---
$ws3 = $wb.Worksheets.Item("All Refactor Variances")
$ws3.Range("I2").Value = 4.879994163690873
$ws3.Range("J2").Value = 3.146150459641948
$ws3.Range("K2").Value = 2.844503589898979
$ws3.Range("L2").Value = 2.231241345921043
$ws3.Range("M2").Value = 2.031125257852269
$ws3.Range("N2").Value = 4.933490050521812
$ws3.Range("O2").Value = 2.898317969195848
$ws3.Range("P2").Value = 2.837310546633725
$ws3.Range("Q2").Value = 2.209999131657667
$ws3.Range("R2").Value = 2.039885453252239
$ws3.Range("I3").Value = 0.2218179165314033
$ws3.Range("J3").Value = 0.143006839074634
$ws3.Range("K3").Value = 0.1292956177226809
$ws3.Range("L3").Value = 0.1014200611782292
$ws3.Range("M3").Value = 0.09232387535692131
$ws3.Range("N3").Value = 0.2349280976438958
$ws3.Range("O3").Value = 0.1380151413902785
$ws3.Range("P3").Value = 0.1351100260301774
$ws3.Range("Q3").Value = 0.1052380538884603
$ws3.Range("R3").Value = 0.0971374025358209
$ws3.Range("I4").Value = 0.2218179165314033
$ws3.Range("J4").Value = 0.3648247556060373
$ws3.Range("K4").Value = 0.4941203733287182
$ws3.Range("L4").Value = 0.5955404345069475
$ws3.Range("M4").Value = 0.6878643098638688
$ws3.Range("N4").Value = 0.2349280976438958
$ws3.Range("O4").Value = 0.3729432390341743
$ws3.Range("P4").Value = 0.5080532650643517
$ws3.Range("Q4").Value = 0.613291318952812
$ws3.Range("R4").Value = 0.7104287214886329
$ws3.Range("I5").Value = 0.3224733618979331
$ws3.Range("J5").Value = 0.2078997805583715
$ws3.Range("K5").Value = 0.1879667484830969
$ws3.Range("L5").Value = 0.1474419587175568
$ws3.Range("M5").Value = 0.1342181503430416
$ws3.Range("N5").Value = 0.3306849660464561
$ws3.Range("O5").Value = 0.1942702163013363
$ws3.Range("P5").Value = 0.1901809737464833
$ws3.Range("Q5").Value = 0.1481331634058156
$ws3.Range("R5").Value = 0.1367306804999087

# --- Sheet: Final Variances ---
$ws4 = $wb.Worksheets.Item("Final Variances")
$ws4.Range("B2").Value = 4.933490050521812
$ws4.Range("C2").Value = 2.898317969195848
$ws4.Range("D2").Value = 2.837310546633725
$ws4.Range("E2").Value = 2.209999131657667
$ws4.Range("F2").Value = 2.039885453252239
$ws4.Range("B3").Value = 0.2349280976438958
$ws4.Range("C3").Value = 0.1380151413902785
$ws4.Range("D3").Value = 0.1351100260301774
$ws4.Range("E3").Value = 0.1052380538884603
$ws4.Range("F3").Value = 0.0971374025358209
$ws4.Range("B4").Value = 0.2349280976438958
$ws4.Range("C4").Value = 0.3729432390341743
$ws4.Range("D4").Value = 0.5080532650643517
$ws4.Range("E4").Value = 0.613291318952812
$ws4.Range("F4").Value = 0.7104287214886329
$ws4.Range("B5").Value = 0.3306849660464561
$ws4.Range("C5").Value = 0.1942702163013363
$ws4.Range("D5").Value = 0.1901809737464833
$ws4.Range("E5").Value = 0.1481331634058156
$ws4.Range("F5").Value = 0.1367306804999087

# --- Sheet: Included and Excluded ---
$ws5 = $wb.Worksheets.Item("Included and Excluded")
$ws5.Range("B2").Value = "[['QEXTRCT', 'QESL', 'QHISPC', 'QEDLESHI', 'PPUNIT', 'QNOHLTH', 'QSERV', 'QFHH', 'PERCAP', 'QNOAUTO', 'QPOVTY', 'MDGRENT', 'QRENTER', 'QFAM', 'QRICH', 'MDHSEVAL', 'QAGEDEP', 'QSSBEN', 'MEDAGE', 'QFEMLBR', 'QFEMALE']]"
